$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the row 2 data: replace Shila Zahra's record with Stella Ireri's record
$ws.Range("D2").Value = "Stella Ireri"
$ws.Range("C2").Value = "stella.ireri@tezzasolutions.com"
$ws.Range("A2").Value = "Stella"
$ws.Range("B2").Value = "Ireri"

# Update the active cell selection to B2 (was E2)
$ws.Range("B2").Select()
